# Add a "pvalue" column (column D) to each worksheet in the workbook.
# Each sheet currently has columns: A=Variable, B=beta, C=se (header row 1,
# data rows 2..N). We append D1="pvalue" (same header style as C1) and the
# pvalue data for rows 2..N.

$wb = $excel.ActiveWorkbook

# Map of worksheet name -> array of p-values for rows 2..N (in row order)
# NOTE: scientific notation (e.g. 1e-06) is not accepted by the PowerShell
# parser here, so all very small p-values are written out in plain decimal
# form instead (value-equivalent once parsed as a double).
$pvaluesBySheet = @{
    "gw" = @(0.1290234611231549,0.3002085914212871,0.2220009610422526,0.1446924608367512,0.7931830271606073,0.01931924298648724,0.7310042812157449,0.2792822672512854,0.7456121300667227,0.851843368544956,0.0608373940281397,0.9762627204220351,0.1847361598184509,0.4675677972722876,0.6226180440369955,0.9277040593429429,0.6670376289183653,0.7306182199557716,0.9061852317161634,0.2623195462746581,0.2454381904917984,0.5973951366002259)
    "shortage" = @(0.0003162805065090288,0.7323874358597695,0.49391430881917,0.003598497022740576,0.003880380508392888,0.02163957395632394,0.000008239817215566001,0.1421789981415316,0.002806223942581709,0.00173225779679854,0.0000002652903284271935,0.07241849061627345,0.9454398639911507,0.3585449778807103,0.009287613265408437)
    "gcpi" = @(0.5010960927306609,0.02538533875937854,0.6953909609254256,0.08810637587206394,0.02735932922402307,0.2452276902447852,0.0001448686781584385,0.1073587853777496,0.7228603072823396,0.7458551455260674,0.637312235973623,0,0.9523798043850926,0.1207433474388797,0.02602729997513364,0.29966872524489,0.002404929925602373,0.6930418784685828,0.9806205159169896,0.9051182364921984,0.3115380699618699,0.0000008614829429820503,0.1294541507853844,0.9721489187917656,0.1340809643009924,0.2828192680549182)
    "cf1" = @(0.001870565130419166,0.0223415990145619,0.07192015529401363,0.1127638405524429,0.0000000000003077538224260934,0.04263410282467461,0.94077935258546,0.5602247091789718,0.03324323423663222,0.003431909069939953,0.000000000004689804100621586,0.9569558931800684,0.6939145477893458,0.1353955003801646)
    "cf10" = @(0.000000000000005773159728050814,0.7143722986940217,0.1070141664693789,0.7146294533782727,0.0003383407290487028,0.2351671397604802,0.1383835634834001,0.9482586225333776,0.8719376760419115)
}

foreach ($sheetName in $pvaluesBySheet.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $values = $pvaluesBySheet[$sheetName]

    # Header cell D1: copy the formatting (bold/border/centered) used by the
    # other header cells (A1:C1), then set its text.
    $ws.Range("C1").Copy()
    $ws.Range("D1").PasteSpecial(-4122)
    $ws.Range("D1").Value = "pvalue"

    # Data cells D2..D(N+1)
    for ($i = 0; $i -lt $values.Count; $i++) {
        $row = $i + 2
        $ws.Cells.Item($row, 4).Value = $values[$i]
    }
}
